$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AF2").Value = "maa://25251 (92.79), ***maa://21730 (25.0), ***maa://39501 (20.83), *maa://36675 (60.0)"
$ws.Range("P3").Value = "maa://21249 (94.32), maa://26254 (96.43)"
$ws.Range("T3").Value = "maa://24617 (89.57), **maa://20790 (44.12), ***maa://37170 (17.19), maa://45854 (100.0)"
$ws.Range("D4").Value = "maa://24632 (93.75), **maa://24303 (33.33), maa://22499 (86.67), maa://22746 (100.0)"
$ws.Range("X4").Value = "**maa://32495 (48.51), ***maa://31785 (22.22), ***maa://36683 (28.26), maa://43217 (86.05)"
$ws.Range("P5").Value = "maa://21919 (96.36), *maa://21281 (80.0)"
$ws.Range("D6").Value = "maa://42407 (95.24)"
$ws.Range("P6").Value = "maa://31836 (92.59), maa://30381 (92.31)"
$ws.Range("X7").Value = "maa://22399 (95.3), *maa://22758 (75.76)"
$ws.Range("AF7").Value = "*maa://26191 (68.67), *maa://36671 (68.0), *maa://42530 (62.5), maa://45272 (100.0)"
$ws.Range("A8").Value = "更新日期：2025.01.30 13:18:46"
$ws.Range("P8").Value = "maa://32931 (83.33), *maa://21916 (60.94), maa://23252 (91.18), maa://37496 (96.77), **maa://22759 (45.45)"
$ws.Range("D9").Value = "maa://22765 (92.39), *maa://21915 (70.37)"
$ws.Range("P9").Value = "maa://22736 (82.47)"
$ws.Range("X9").Value = "maa://26223 (97.71)"
$ws.Range("D10").Value = "***maa://25695 (18.82), **maa://32237 (41.3), ***maa://34206 (20.0), ***maa://39951 (15.22), ***maa://39243 (28.57), *maa://45271 (55.0)"
$ws.Range("P10").Value = "maa://28977 (90.36), maa://36669 (88.1), *maa://23264 (61.82)"
$ws.Range("X10").Value = "maa://22301 (97.74), maa://22726 (100.0), maa://45828 (100.0)"
$ws.Range("T11").Value = "maa://22747 (92.9), maa://22501 (97.59), *maa://45521 (66.67)"
$ws.Range("X12").Value = "maa://22753 (91.28), *maa://21485 (76.26), maa://37962 (88.57)"
$ws.Range("D13").Value = "maa://24999 (91.89), maa://36673 (92.96), maa://25001 (85.51)"
$ws.Range("H13").Value = "*maa://21248 (73.71), **maa://22728 (47.73)"
$ws.Range("P13").Value = "maa://22676 (92.62), *maa://22583 (74.24), *maa://22500 (58.7)"
$ws.Range("AF13").Value = "**maa://22737 (33.57), maa://39883 (92.19), *maa://39885 (57.14)"
$ws.Range("P14").Value = "maa://23250 (98.69), maa://20107 (87.1), maa://22772 (100.0), **maa://22745 (50.0)"
$ws.Range("AB14").Value = "maa://22764 (97.1)"
$ws.Range("D15").Value = "*maa://22743 (77.83), maa://22734 (84.03), *maa://30808 (64.18), **maa://36048 (42.86), maa://45058 (100.0)"
$ws.Range("H15").Value = "maa://24304 (88.1), maa://21478 (91.67)"
$ws.Range("AF15").Value = "maa://21364 (81.37), *maa://36666 (78.12), *maa://22766 (68.7)"
$ws.Range("D16").Value = "maa://21441 (96.4), maa://36679 (93.75), maa://37650 (97.06)"
$ws.Range("T16").Value = "maa://22729 (94.81), *maa://28648 (68.75), maa://36674 (80.43)"
$ws.Range("X16").Value = "maa://28501 (97.98), maa://28051 (96.0)"
$ws.Range("AB16").Value = "maa://26228 (95.83)"
$ws.Range("AF16").Value = "*maa://23911 (64.76), maa://27755 (93.1)"
$ws.Range("H17").Value = "maa://22430 (88.66), maa://39599 (85.11)"
$ws.Range("H18").Value = "maa://24421 (89.2)"
$ws.Range("L18").Value = "maa://22466 (89.47), *maa://22732 (51.16)"
$ws.Range("T18").Value = "maa://24385 (97.22)"
$ws.Range("D20").Value = "maa://21432 (89.94), maa://25198 (93.4), *maa://20795 (51.56), maa://36680 (93.75)"
$ws.Range("H20").Value = "maa://22864 (89.81)"
$ws.Range("L20").Value = "maa://41331 (85.07)"
$ws.Range("AB21").Value = "maa://21443 (80.71), ***maa://23820 (29.31)"
$ws.Range("AF21").Value = "maa://22524 (94.37), *maa://22432 (78.46)"
$ws.Range("L22").Value = "maa://27127 (81.58), *maa://22751 (72.06)"
$ws.Range("X22").Value = "maa://21282 (98.57), *maa://37649 (67.86)"
$ws.Range("D23").Value = "***maa://28036 (28.57), *maa://41753 (57.14)"
$ws.Range("L23").Value = "maa://39756 (95.19), maa://39875 (94.12)"
$ws.Range("D24").Value = "*maa://24368 (77.96)"
$ws.Range("X24").Value = "maa://29988 (85.02), maa://23504 (93.15), **maa://22892 (40.14), *maa://25141 (76.74), *maa://36663 (78.08), ***maa://22815 (23.08)"
$ws.Range("AF24").Value = "maa://22523 (85.86), maa://36672 (80.36), maa://29910 (92.86), **maa://21440 (34.55), maa://45831 (100.0)"
$ws.Range("D25").Value = "maa://29753 (95.06)"
$ws.Range("H25").Value = "*maa://29063 (74.52), *maa://25311 (73.53), ***maa://22725 (4.84), *maa://45047 (66.67)"
$ws.Range("T25").Value = "maa://20109 (92.49), maa://22545 (100.0), maa://42915 (100.0)"
$ws.Range("X25").Value = "*maa://29890 (78.26)"
$ws.Range("AB25").Value = "maa://31215 (87.16), *maa://24516 (80.0), maa://26001 (87.5)"
$ws.Range("AF25").Value = "maa://20108 (96.3), maa://24621 (96.75), maa://36676 (96.88), maa://22771 (85.71), *maa://37772 (66.67)"
$ws.Range("H26").Value = "maa://24913 (91.86)"
$ws.Range("AB26").Value = "maa://42235 (94.51)"
$ws.Range("T28").Value = "maa://23263 (95.19), *maa://29765 (63.41)"
$ws.Range("X28").Value = "maa://39929 (90.41), maa://41749 (91.89), ***maa://39723 (13.89)"
$ws.Range("AF28").Value = "maa://36660 (92.51), *maa://36701 (65.52)"
$ws.Range("H29").Value = "*maa://25175 (66.67)"
$ws.Range("AF29").Value = "*maa://24080 (68.77), maa://42865 (82.0), ***maa://34960 (8.33)"
$ws.Range("AB30").Value = "maa://42979 (96.43), maa://45822 (100.0), maa://45045 (100.0)"
$ws.Range("L31").Value = "maa://35926 (93.53), maa://36258 (85.32), *maa://43904 (72.73)"
$ws.Range("L32").Value = "maa://28065 (95.24)"
$ws.Range("T32").Value = "maa://42859 (96.12), maa://41108 (88.0), maa://41238 (96.97), maa://45523 (100.0)"
$ws.Range("T34").Value = "maa://24526 (93.28)"
$ws.Range("L35").Value = "maa://41296 (95.92)"
$ws.Range("P37").Value = "maa://21280 (89.2), *maa://21239 (66.67)"
$ws.Range("T37").Value = "**maa://39354 (45.45)"
$ws.Range("P38").Value = "*maa://24383 (69.61)"
$ws.Range("AF38").Value = "maa://36697 (87.06)"
$ws.Range("H39").Value = "maa://36670 (88.17), maa://25199 (84.82), maa://30434 (90.79), ***maa://25036 (16.0), *maa://45059 (75.0), *maa://44165 (66.67)"
$ws.Range("T39").Value = "maa://45788 (80.85), maa://45790 (88.89)"
$ws.Range("P40").Value = "maa://23278 (95.77), maa://21386 (95.74), maa://36664 (90.91), maa://45550 (100.0)"
$ws.Range("H44").Value = "maa://29768 (97.94), maa://27728 (96.04)"
$ws.Range("H45").Value = "maa://21229 (84.49), maa://30807 (95.59), *maa://22767 (55.0), ***maa://20796 (13.79), maa://42459 (84.21)"
$ws.Range("H46").Value = "maa://35931 (92.41), maa://43901 (88.24)"
$ws.Range("H53").Value = "maa://32534 (93.75), **maa://32434 (33.33)"
$ws.Range("H55").Value = "maa://32532 (92.45)"
$ws.Range("H60").Value = "*maa://40438 (67.27)"
$ws.Range("H62").Value = "maa://42981 (97.22), maa://43903 (100.0)"
